$d = $word.ActiveDocument

# 1) Add a new row to the first table with Luca Cappon's name and email,
#    matching the style of the existing row.
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Luca Cappon (2037394)"
$newRow.Cells.Item(2).Range.Text = "h.l.cappon@student.tue.nl"

# 2) Merge the two runs "Changing font, spacing, or document " and
#    "layout is not allowed." into a single run of text.
$d.Content.Find.Execute("Changing font, spacing, or document layout is not allowed.", $false, $false, $false, $false, $false, $true, 1, $false, "Changing font, spacing, or document layout is not allowed.", 2) | Out-Null

# 3) Mark the picture's run as NoProof (adds <w:noProof/> to its run
#    properties) for the screenshot inline image.
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = $true
